$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new sprint-backlog row above what used to be row 5 ("Rank close
# parking spots by driving distance"), shifting it and everything below
# down by one.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row with the new backlog item.
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Implement coordinates class"
$ws.Range("C5").Value = 12

# Rename the Estimated Time column header to include units, and open the
# sprint backlog for the sprint meeting with updated hour estimates.
$ws.Range("C1").Value = "Estimated Time (hr)"

$ws.Range("C3").Value = 15
$ws.Range("C4").Value = 1
$ws.Range("C6").Value = 10
$ws.Range("C7").Value = 40
$ws.Range("C9").Value = 5
$ws.Range("C10").Value = 1
$ws.Range("C11").Value = 20
$ws.Range("C13").Value = 50
$ws.Range("C14").Value = 3

$ws.Range("D5").Select()
